$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.783.63"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").Value = "1.625.85"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.06"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5108"
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2582"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06377"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  -2.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07771"
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.249"
$ws.Range("E12").Value = "  -1.01%  "
$ws.Range("D13").Value = "1.625.22"
$ws.Range("E13").Value = "  -0.97%  "
$ws.Range("D14").Value = "1.850.17"
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5559"
$ws.Range("E15").Value = "  +1.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.48"
$ws.Range("E16").Value = "  -1.69%  "
$ws.Range("D17").Value = "0.0₅7530"
$ws.Range("E17").Value = "  -3.24%  "
$ws.Range("D18").Value = "25.779.42"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.82"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.334"
$ws.Range("E21").Value = "  -3.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.779"
$ws.Range("E22").Value = "  -2.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.992"
$ws.Range("E23").Value = "  -1.94%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.819"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1289"
$ws.Range("E26").Value = "  +4.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "141.31"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.735"
$ws.Range("E28").Value = "  -2.07%  "
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("E30").Value = "  -0.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04880"
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.295"
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.186"
$ws.Range("E33").Value = "  -1.43%  "
$ws.Range("E34").Value = "  +1.00%  "
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8943"
$ws.Range("E36").Value = "  -2.69%  "
$ws.Range("D37").Value = "1.134.01"
$ws.Range("E37").Value = "  +0.92%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.539"
$ws.Range("E38").Value = "  -1.21%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5486"
$ws.Range("E39").Value = "  -2.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01559"
$ws.Range("E40").Value = "  -1.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9952"
$ws.Range("E41").Value = "  -0.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.590"
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7942"
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.25"
$ws.Range("E44").Value = "  -2.50%  "
$ws.Range("D45").Value = "1.774.91"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("D46").Value = "0.0₈112"
$ws.Range("E46").Value = "  -7.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4423"
$ws.Range("E47").Value = "  -2.49%  "
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05065"
$ws.Range("E49").Value = "  -2.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.549"
$ws.Range("E50").Value = "  +0.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.004"
$ws.Range("E51").Value = "  -0.34%  "
